$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '328.37'
Set-TextValue 'E2' '1.55%'
Set-TextValue 'D3' '41.58'
Set-TextValue 'E3' '5.20%'
Set-TextValue 'D4' '5.615'
Set-TextValue 'E4' '-4.33%'
Set-TextValue 'D5' '0.08174'
Set-TextValue 'E5' '1.84%'
Set-TextValue 'D6' '2.023'
Set-TextValue 'E6' '1.25%'
Set-TextValue 'D7' '8.734'
Set-TextValue 'E7' '1.36%'
Set-TextValue 'D8' '4.522'
Set-TextValue 'E8' '-1.07%'
Set-TextValue 'D9' '2.983'
Set-TextValue 'E9' '1.22%'
Set-TextValue 'D10' '0.9218'
Set-TextValue 'E10' '-0.76%'
Set-TextValue 'D11' '0.1277'
Set-TextValue 'E11' '0.63%'
Set-TextValue 'D12' '0.1960'
Set-TextValue 'E12' '0.65%'
Set-TextValue 'D13' '0.09415'
Set-TextValue 'E13' '3.07%'
Set-TextValue 'D14' '0.03808'
Set-TextValue 'E14' '4.98%'
Set-TextValue 'E15' '0.91%'
Set-TextValue 'D16' '0.001307'
Set-TextValue 'E16' '1.16%'
Set-TextValue 'D17' '0.006282'
Set-TextValue 'E17' '-0.94%'
Set-TextValue 'D19' '3.445'
Set-TextValue 'E19' '2.76%'
Set-TextValue 'E20' '-1.22%'
Set-TextValue 'D21' '8.317'
Set-TextValue 'E21' '-4.42%'
Set-TextValue 'D22' '0.1387'
Set-TextValue 'E22' '1.18%'
Set-TextValue 'D23' '0.2412'
Set-TextValue 'E23' '-1.44%'
Set-TextValue 'E24' '-0.06%'
Set-TextValue 'D25' '0.001259'
Set-TextValue 'E25' '-0.33%'
Set-TextValue 'D26' '0.004340'
Set-TextValue 'E26' '-1.04%'
Set-TextValue 'D27' '0.0001182'
Set-TextValue 'E27' '2.63%'
Set-TextValue 'D39' '0.02775'
Set-TextValue 'E39' '10.85%'
Set-TextValue 'D40' '0.05426'
Set-TextValue 'E40' '3.97%'
Set-TextValue 'D41' '0.007664'
Set-TextValue 'E41' '2.65%'
Set-TextValue 'D42' '0.1419'
Set-TextValue 'E42' '1.17%'
Set-TextValue 'D43' '0.008979'
Set-TextValue 'E43' '-6.68%'
Set-TextValue 'D44' '0.002134'
Set-TextValue 'E44' '0.67%'
Set-TextValue 'D45' '0.01169'
Set-TextValue 'E45' '6.03%'
Set-TextValue 'D46' '0.00006657'
Set-TextValue 'E46' '-1.49%'
Set-TextValue 'E47' '0.03%'
Set-TextValue 'D48' '0.003201'
Set-TextValue 'E48' '6.53%'
Set-TextValue 'D49' '0.002282'
Set-TextValue 'E49' '-0.49%'
Set-TextValue 'D50' '0.00002103'
Set-TextValue 'E50' '0.03%'
Set-TextValue 'E51' '0.03%'
